$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 7.333333333333334
    3  = 7.333333333333334
    4  = -3.666666666666667
    5  = -7.333333333333334
    6  = -7.333333333333334
    7  = 3.666666666666667
    8  = 7.333333333333334
    9  = 11
    10 = 7.333333333333334
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
